# Generate Report for Handoff
#
# - Overview sheet: bump "Latest HO Xliff Generate Date" (col G) for the
#   rows that were just (re)handed off.
# - zh-cn / de-de sheets: bump "Latest Handoff Datetime" (col H) for the
#   same rows, and mark their "Priority" (col E) as "ht" (handed off).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 11, 12, 13, 14)

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-27 02:21:49"
}

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-27 02:21:44"
}

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-27 02:21:49"
}
